# //VJ 170210 FIX BUG: - stemflow causes MB error and drowned interception also.
# - channel overlfow whlevel < 0, no movement
# - channelDX = DX
# - Clean up some code
#
# Update the flood-channel input parameters on Sheet1 (second scenario block,
# rows 35-46) and refresh the view/selection state to match the saved
# workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$null = $ws.Activate()

# --- Input values changed for the second scenario (rows 35-46) ---
# cd (C36): 19 -> 5
$ws.Range("C36").Value = 5
# V   (C39): 0.7 -> 0.2
$ws.Range("C39").Value = 0.2
# V   (C40): 1.8 -> 0.1
$ws.Range("C40").Value = 0.1

# --- Column B width tweak (bestFit re-measured by newer Excel build) ---
$ws.Columns("B").ColumnWidth = 8.8

# --- View state: scroll position + active selection moved to D45 ---
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("D45").Select()
